# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
#
# Five fixtures had their data (everything except the running "id" in
# column A) crossed between two neighbouring rows. This swaps each pair
# of rows' B:AD contents back into the correct row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($RowA, $RowB)

    $addrA = "B" + $RowA + ":AD" + $RowA
    $addrB = "B" + $RowB + ":AD" + $RowB

    $rangeA = $ws.Range($addrA)
    $rangeB = $ws.Range($addrB)

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-RowData 11 12
Swap-RowData 19 22
Swap-RowData 27 28
Swap-RowData 43 44
Swap-RowData 47 48
